# Insert a new data row before the current row 253 (Vega Monumental Concepción /
# Brócoli sheet), shifting all subsequent rows down by one, and populate the
# newly inserted row with a new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 253; this shifts old rows 253..311 down to 254..312
# and copies formatting (e.g. the date number format on column D) from the row
# above, matching how Excel behaves when inserting rows.
$ws.Rows("253").Insert()

# Populate the new row 253 with the new record's values.
$ws.Range("A253").Value = 11
$ws.Range("B253").Value = "Vega Monumental Concepción"
$ws.Range("C253").Value = "Bíobío"
$ws.Range("D253").Value = 44754
$ws.Range("E253").Value = 8
$ws.Range("F253").Value = 100112023
$ws.Range("G253").Value = "Brócoli"
$ws.Range("H253").Value = "Sin especificar"
$ws.Range("I253").Value = "Primera"
$ws.Range("J253").Value = 2000
$ws.Range("K253").Value = 1000
$ws.Range("L253").Value = 1100
$ws.Range("M253").Value = 1050
$ws.Range("N253").Value = "`$/unidad"
$ws.Range("O253").Value = "Región Metropolitana"
$ws.Range("P253").Value = 1050
$ws.Range("Q253").Value = 1
$ws.Range("R253").Value = "Hortaliza"
